# Generate Report for Handoff
#
# The localization CI run moved each language from "handed back" to ready for
# a fresh handoff: the status text flips to "Ready for handoff" and the
# handoff timestamps are bumped. Because the new status string is shorter
# than the old one, the (auto-sized) status columns also shrink.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet: E2/F2 are the per-language status, G2 the handoff date
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-17 15:01:39"

# --- zh-cn sheet: C2 status, H2 latest handoff datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-17 15:01:33"

# --- de-de sheet: C2 status, H2 latest handoff datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-17 15:01:39"

# --- Re-fit the status columns now that "Ready for handoff" is shorter than
# "Handed back: in sync with en-US"
$newStatusColumnWidth = 16.38265482584637
$overview.Columns.Item(5).ColumnWidth = $newStatusColumnWidth
$overview.Columns.Item(6).ColumnWidth = $newStatusColumnWidth
$zhcn.Columns.Item(3).ColumnWidth = $newStatusColumnWidth
$dede.Columns.Item(3).ColumnWidth = $newStatusColumnWidth
